{"js": "// Replace the date line and the twenty-five \"A\u00d7B=C\" answer cells with\n// their updated values, per the authored diff.\nconst replacements = [\n  [\"2025-09-26 Friday\", \"2025-09-27 Saturday\"],\n  [\"237\u00d77=1659\", \"481\u00d77=3367\"],\n  [\"621\u00d75=3105\", \"982\u00d79=8838\"],\n  [\"659\u00d76=3954\", \"241\u00d74=964\"],\n  [\"523\u00d74=2092\", \"959\u00d73=2877\"],\n  [\"206\u00d79=1854\", \"428\u00d79=3852\"],\n  [\"795\u00d74=3180\", \"368\u00d78=2944\"],\n  [\"773\u00d75=3865\", \"548\u00d77=3836\"],\n  [\"969\u00d72=1938\", \"738\u00d76=4428\"],\n  [\"198\u00d75=990\", \"644\u00d79=5796\"],\n  [\"646\u00d75=3230\", \"505\u00d75=2525\"],\n  [\"219\u00d73=657\", \"530\u00d72=1060\"],\n  [\"106\u00d77=742\", \"585\u00d78=4680\"],\n  [\"548\u00d78=4384\", \"313\u00d79=2817\"],\n  [\"826\u00d74=3304\", \"198\u00d72=396\"],\n  [\"526\u00d77=3682\", \"989\u00d76=5934\"],\n  [\"832\u00d78=6656\", \"240\u00d73=720\"],\n  [\"844\u00d73=2532\", \"386\u00d72=772\"],\n  [\"607\u00d76=3642\", \"775\u00d74=3100\"],\n  [\"639\u00d74=2556\", \"838\u00d77=5866\"],\n  [\"880\u00d76=5280\", \"549\u00d73=1647\"],\n  [\"865\u00d73=2595\", \"678\u00d76=4068\"],\n  [\"597\u00d77=4179\", \"431\u00d76=2586\"],\n  [\"183\u00d76=1098\", \"903\u00d78=7224\"],\n  [\"494\u00d79=4446\", \"569\u00d77=3983\"],\n  [\"455\u00d76=2730\", \"588\u00d72=1176\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the twenty-five \"A\u00d7B=C\" answer cells with\n# their updated values, per the authored diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-09-26 Friday\", \"2025-09-27 Saturday\"),\n  @(\"237\u00d77=1659\", \"481\u00d77=3367\"),\n  @(\"621\u00d75=3105\", \"982\u00d79=8838\"),\n  @(\"659\u00d76=3954\", \"241\u00d74=964\"),\n  @(\"523\u00d74=2092\", \"959\u00d73=2877\"),\n  @(\"206\u00d79=1854\", \"428\u00d79=3852\"),\n  @(\"795\u00d74=3180\", \"368\u00d78=2944\"),\n  @(\"773\u00d75=3865\", \"548\u00d77=3836\"),\n  @(\"969\u00d72=1938\", \"738\u00d76=4428\"),\n  @(\"198\u00d75=990\", \"644\u00d79=5796\"),\n  @(\"646\u00d75=3230\", \"505\u00d75=2525\"),\n  @(\"219\u00d73=657\", \"530\u00d72=1060\"),\n  @(\"106\u00d77=742\", \"585\u00d78=4680\"),\n  @(\"548\u00d78=4384\", \"313\u00d79=2817\"),\n  @(\"826\u00d74=3304\", \"198\u00d72=396\"),\n  @(\"526\u00d77=3682\", \"989\u00d76=5934\"),\n  @(\"832\u00d78=6656\", \"240\u00d73=720\"),\n  @(\"844\u00d73=2532\", \"386\u00d72=772\"),\n  @(\"607\u00d76=3642\", \"775\u00d74=3100\"),\n  @(\"639\u00d74=2556\", \"838\u00d77=5866\"),\n  @(\"880\u00d76=5280\", \"549\u00d73=1647\"),\n  @(\"865\u00d73=2595\", \"678\u00d76=4068\"),\n  @(\"597\u00d77=4179\", \"431\u00d76=2586\"),\n  @(\"183\u00d76=1098\", \"903\u00d78=7224\"),\n  @(\"494\u00d79=4446\", \"569\u00d77=3983\"),\n  @(\"455\u00d76=2730\", \"588\u00d72=1176\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
